$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rows = @(
  @(2, 46048.01041666666, 658),
  @(3, 46048.02083333334, 671),
  @(4, 46048.03125, 714),
  @(5, 46048.04166666666, 733),
  @(6, 46048.05208333334, 712),
  @(7, 46048.0625, 717),
  @(8, 46048.07291666666, 709),
  @(9, 46048.08333333334, 703),
  @(10, 46048.09375, 717),
  @(11, 46048.10416666666, 737),
  @(12, 46048.11458333334, 784),
  @(13, 46048.125, 833),
  @(14, 46048.13541666666, 867),
  @(15, 46048.14583333334, 899),
  @(16, 46048.15625, 920),
  @(17, 46048.16666666666, 955),
  @(18, 46048.17708333334, 971),
  @(19, 46048.1875, 980),
  @(20, 46048.19791666666, 972),
  @(21, 46048.20833333334, 993),
  @(22, 46048.21875, 1053),
  @(23, 46048.22916666666, 1097),
  @(24, 46048.23958333334, 1191),
  @(25, 46048.25, 1276),
  @(26, 46048.26041666666, 1337),
  @(27, 46048.27083333334, 1457),
  @(28, 46048.28125, 1527),
  @(29, 46048.29166666666, 1554),
  @(30, 46048.30208333334, 1601),
  @(31, 46048.3125, 1644),
  @(32, 46048.32291666666, 1691),
  @(33, 46048.33333333334, 0),
  @(34, 46048.34375, 1631),
  @(35, 46048.35416666666, 1645),
  @(36, 46048.36458333334, 1643),
  @(37, 46048.375, 1637),
  @(38, 46048.38541666666, 1569),
  @(39, 46048.39583333334, 1499),
  @(40, 46048.40625, 1468),
  @(41, 46048.41666666666, 1498),
  @(42, 46048.42708333334, 1495),
  @(43, 46048.4375, 0),
  @(44, 46048.44791666666, 0),
  @(45, 46048.45833333334, 0),
  @(46, 46048.46875, 0),
  @(47, 46048.47916666666, 0),
  @(48, 46048.48958333334, 0),
  @(49, 46048.5, 0),
  @(50, 46048.51041666666, 0),
  @(51, 46048.52083333334, 0),
  @(52, 46048.53125, 0),
  @(53, 46048.54166666666, 0),
  @(54, 46048.55208333334, 0),
  @(55, 46048.5625, 0),
  @(56, 46048.57291666666, 0),
  @(57, 46048.58333333334, 0),
  @(58, 46048.59375, 0),
  @(59, 46048.60416666666, 0),
  @(60, 46048.61458333334, 0),
  @(61, 46048.625, 0),
  @(62, 46048.63541666666, 0),
  @(63, 46048.64583333334, 0),
  @(64, 46048.65625, 0),
  @(65, 46048.66666666666, 0),
  @(66, 46048.67708333334, 0),
  @(67, 46048.6875, 0),
  @(68, 46048.69791666666, 0),
  @(69, 46048.70833333334, 0),
  @(70, 46048.71875, 0),
  @(71, 46048.72916666666, 0),
  @(72, 46048.73958333334, 0),
  @(73, 46048.75, 0),
  @(74, 46048.76041666666, 0),
  @(75, 46048.77083333334, 0),
  @(76, 46048.78125, 0),
  @(77, 46048.79166666666, 0),
  @(78, 46048.80208333334, 0),
  @(79, 46048.8125, 0),
  @(80, 46048.82291666666, 0),
  @(81, 46048.83333333334, 0),
  @(82, 46048.84375, 0),
  @(83, 46048.85416666666, 0),
  @(84, 46048.86458333334, 0),
  @(85, 46048.875, 0),
  @(86, 46048.88541666666, 0),
  @(87, 46048.89583333334, 0),
  @(88, 46048.90625, 0),
  @(89, 46048.91666666666, 0),
  @(90, 46048.92708333334, 0),
  @(91, 46048.9375, 0),
  @(92, 46048.94791666666, 0),
  @(93, 46048.95833333334, 0),
  @(94, 46048.96875, 0),
  @(95, 46048.97916666666, 0),
  @(96, 46048.98958333334, 0),
  @(97, 46049, 0)
)
foreach ($row in $rows) {
    $r = $row[0]
    $a = $row[1]
    $b = $row[2]
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
}
